$data = @{
    2 = @(8, "{""L2"":4,""L1"":4}", "{""L2"":0.5,""L1"":0.5}", 0.5, "50.00", "L2")
    3 = $null
    4 = @(15, "{""L2"":7,""L3"":7,""L1"":1}", "{""L2"":0.4666666667,""L3"":0.4666666667,""L1"":0.0666666667}", 0.4666666666666667, "46.67", "L3")
    5 = $null
    6 = @(8, "{""L2"":8}", "{""L2"":1.0}", 1, "100.00", "L2")
    7 = @(8, "{""L1"":7,""L2"":1}", "{""L1"":0.875,""L2"":0.125}", 0.875, "87.50", "L2")
    8 = @(8, "{""L2"":8}", "{""L2"":1.0}", 1, "100.00", "L2")
    9 = @(8, "{""L2"":8}", "{""L2"":1.0}", 1, "100.00", "L2")
    10 = @(8, "{""L2"":8}", "{""L2"":1.0}", 1, "100.00", "L2")
    11 = @(8, "{""L2"":8}", "{""L2"":1.0}", 1, "100.00", "L2")
    12 = @(8, "{""L1"":7,""L2"":1}", "{""L1"":0.875,""L2"":0.125}", 0.875, "87.50", "L2")
    13 = @(8, "{""L2"":7,""L3"":1}", "{""L2"":0.875,""L3"":0.125}", 0.875, "87.50", "L3")
    14 = $null
    15 = $null
    16 = @(8, "{""L3"":8}", "{""L3"":1.0}", 1, "100.00", "L3")
    17 = $null
    18 = @(8, "{""L2"":6,""L1"":2}", "{""L2"":0.75,""L1"":0.25}", 0.75, "75.00", "L2")
    19 = $null
    20 = @(8, "{""L3"":7,""L2"":1}", "{""L3"":0.875,""L2"":0.125}", 0.875, "87.50", "L3")
    21 = @(8, "{""L3"":6,""L2"":2}", "{""L3"":0.75,""L2"":0.25}", 0.75, "75.00", "L3")
    22 = @(8, "{""L2"":7,""L1"":1}", "{""L2"":0.875,""L1"":0.125}", 0.875, "87.50", "L2")
    23 = @(8, "{""L3"":7,""L2"":1}", "{""L3"":0.875,""L2"":0.125}", 0.875, "87.50", "L3")
    24 = @(8, "{""L3"":7,""L2"":1}", "{""L3"":0.875,""L2"":0.125}", 0.875, "87.50", "L3")
    25 = @(8, "{""L3"":8}", "{""L3"":1.0}", 1, "100.00", "L3")
    26 = @(8, "{""L3"":8}", "{""L3"":1.0}", 1, "100.00", "L3")
    27 = $null
    28 = $null
    29 = $null
    30 = $null
    31 = $null
    32 = $null
    33 = $null
    34 = $null
    35 = $null
    36 = $null
    37 = $null
    38 = $null
    39 = $null
    40 = $null
    41 = $null
    42 = $null
    43 = $null
    44 = @(8, "{""L2"":5,""L1"":3}", "{""L2"":0.625,""L1"":0.375}", 0.625, "62.50", "L2")
    45 = $null
    46 = $null
    47 = $null
    48 = $null
    49 = $null
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: column F's header changes from "recommended_level" to
#     "frequency"; five new columns are inserted (G:K) carrying the
#     frequency-analysis headers, with K1 taking over the old
#     "recommended_level" title. Copy F1's formatting (bold, centered,
#     bordered) onto the new header cells before writing their text.
$ws.Range("F1").Copy()
$ws.Range("G1:K1").PasteSpecial(-4122)
$ws.Range("F1").Value = "frequency"
$ws.Range("G1").Value = "frequency_occurrence"
$ws.Range("H1").Value = "frequency_occurrence_probab"
$ws.Range("I1").Value = "max_probab"
$ws.Range("J1").Value = "max_probab_percentage"
$ws.Range("K1").Value = "recommended_level"

# --- Data rows: the old column F (recommended_level) moves to the new
#     column K, and F:J get populated with the frequency-analysis figures
#     (frequency count, per-level occurrence counts/probabilities as JSON,
#     the winning probability, and its percentage string).
for ($row = 2; $row -le 49; $row++) {
    # Value2 (unlike Value) returns a clean $null for a blank source cell.
    $oldRecommended = $ws.Cells.Item($row, 6).Value2
    $ws.Cells.Item($row, 11).Value = $oldRecommended

    $vals = $data[$row]
    if ($vals -ne $null) {
        $ws.Cells.Item($row, 6).Value = $vals[0]
        $ws.Cells.Item($row, 7).Value = $vals[1]
        $ws.Cells.Item($row, 8).Value = $vals[2]
        $ws.Cells.Item($row, 9).Value = $vals[3]

        # Column J holds a numeric-looking percentage string ("50.00",
        # "100.00", ...) that Excel would otherwise auto-coerce to a
        # number; force text with a quote-prefix, then drop the style
        # back to Normal so no stray number-format sticks to the cell.
        $ws.Cells.Item($row, 10).Value = "'" + $vals[4]
        $ws.Cells.Item($row, 10).Style = "Normal"
    } else {
        $ws.Cells.Item($row, 6).Value = $oldRecommended
    }
}
